# Refresh the cryptocurrency price/volume snapshot (scheduled GitHub Actions update).
# Only the Price (D) and Volume(1h) (E) columns move for most rows; two rows (44/45)
# also swap which coin (B/C) occupies that rank.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.614.41'
$ws.Range('E2').Value = '  +0.50%  '

$ws.Range('D3').Value = '2.818.86'
$ws.Range('E3').Value = '  +1.61%  '

$ws.Range('E4').Value = '  +0.09%  '

$ws.Range('D5').Value = '''349.93'
$ws.Range('E5').Value = '  -0.82%  '

$ws.Range('D6').Value = '''112.22'
$ws.Range('E6').Value = '  +3.84%  '

$ws.Range('D7').Value = '''0.556'
$ws.Range('E7').Value = '  +1.24%  '

$ws.Range('E8').Value = '  +0.05%  '

$ws.Range('D9').Value = '''0.615'
$ws.Range('E9').Value = '  +4.85%  '

$ws.Range('D10').Value = '''39.91'
$ws.Range('E10').Value = '  +0.80%  '

$ws.Range('E11').Value = '  -0.74%  '

$ws.Range('D12').Value = '''0.0845'
$ws.Range('E12').Value = '  +1.22%  '

$ws.Range('D13').Value = '''19.86'
$ws.Range('E13').Value = '  -0.88%  '

$ws.Range('D14').Value = '''7.73'
$ws.Range('E14').Value = '  +2.36%  '

$ws.Range('D15').Value = '3.270.70'
$ws.Range('E15').Value = '  +1.91%  '

$ws.Range('E16').Value = '  +5.47%  '

$ws.Range('D17').Value = '2.816.69'
$ws.Range('E17').Value = '  +1.46%  '

$ws.Range('D18').Value = '51.700.65'
$ws.Range('E18').Value = '  +0.70%  '

$ws.Range('D19').Value = '''3.44'
$ws.Range('E19').Value = '  +10.95%  '

$ws.Range('D20').Value = '''7.57'
$ws.Range('E20').Value = '  -0.28%  '

$ws.Range('D21').Value = '''13.33'
$ws.Range('E21').Value = '  +1.57%  '

$ws.Range('D22').Value = '0.0₃0969'
$ws.Range('E22').Value = '  +0.73%  '

$ws.Range('D23').Value = '''70.30'
$ws.Range('E23').Value = '  +0.68%  '

$ws.Range('D24').Value = '''267.53'
$ws.Range('E24').Value = '  +0.80%  '

$ws.Range('D25').Value = '''2.74'
$ws.Range('E25').Value = '  +1.30%  '

$ws.Range('D26').Value = '''1.00'
$ws.Range('E26').Value = '  +0.03%  '

$ws.Range('D27').Value = '''26.13'
$ws.Range('E27').Value = '  +0.45%  '

$ws.Range('D28').Value = '''0.162'
$ws.Range('E28').Value = '  +0.23%  '

$ws.Range('D29').Value = '''10.53'
$ws.Range('E29').Value = '  +3.09%  '

$ws.Range('D30').Value = '''38.32'
$ws.Range('E30').Value = '  +5.43%  '

$ws.Range('E31').Value = '  +2.83%  '

$ws.Range('D32').Value = '''6.29'
$ws.Range('E32').Value = '  +2.09%  '

$ws.Range('D33').Value = '''52.76'
$ws.Range('E33').Value = '  +1.62%  '

$ws.Range('D34').Value = '''0.0892'
$ws.Range('E34').Value = '  +7.98%  '

$ws.Range('D35').Value = '''0.0447'
$ws.Range('E35').Value = '  -1.27%  '

$ws.Range('D36').Value = '''5.59'
$ws.Range('E36').Value = '  +1.09%  '

$ws.Range('E37').Value = '  +0.10%  '

$ws.Range('D38').Value = '''18.86'
$ws.Range('E38').Value = '  +2.45%  '

$ws.Range('D39').Value = '''3.20'
$ws.Range('E39').Value = '  +1.71%  '

$ws.Range('D40').Value = '''2.00'
$ws.Range('E40').Value = '  +2.09%  '

$ws.Range('D41').Value = '''0.115'
$ws.Range('E41').Value = '  +1.34%  '

$ws.Range('D42').Value = '''2.50'
$ws.Range('E42').Value = '  -1.02%  '

$ws.Range('D43').Value = '''122.67'
$ws.Range('E43').Value = '  +1.90%  '

$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').Value = '''2.21'
$ws.Range('E44').Value = '  +1.06%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '''22.03'
$ws.Range('E45').Value = '  -0.08%  '

$ws.Range('E46').Value = '  +8.59%  '

$ws.Range('D47').Value = '2.173.07'
$ws.Range('E47').Value = '  +3.41%  '

$ws.Range('D48').Value = '''3.47'
$ws.Range('E48').Value = '  +6.75%  '

$ws.Range('D49').Value = '''0.250'
$ws.Range('E49').Value = '  +25.43%  '

$ws.Range('D50').Value = '''0.944'
$ws.Range('E50').Value = '  +4.46%  '

$ws.Range('D51').Value = '''5.50'
$ws.Range('E51').Value = '  +1.72%  '
